# Apply updated cryptocurrency price/volume data to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, never letting Excel
# auto-convert numeric-looking strings ("143.32", "0.0000141", ...) into
# real numbers. A leading apostrophe forces text entry (like a user typing
# '143.32 into the cell); resetting the Style back to 'Normal' afterwards
# drops the transient quote-prefix style so the cell keeps its original,
# unstyled General formatting.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = 'Normal'
}


# Row 2
Set-TextValue $ws.Range('D2') '60.917.35'
$ws.Range('E2').Value = '  +2.56%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.614.07'
$ws.Range('E3').Value = '  +0.93%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
Set-TextValue $ws.Range('D5') '574.49'
$ws.Range('E5').Value = '  +0.56%  '

# Row 6
Set-TextValue $ws.Range('D6') '143.32'
$ws.Range('E6').Value = '  -0.46%  '

# Row 7
$ws.Range('E7').Value = '  -0.24%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.602'
$ws.Range('E8').Value = '  +0.76%  '

# Row 9
Set-TextValue $ws.Range('D9') '2.639.51'
$ws.Range('E9').Value = '  +1.54%  '

# Row 10
$ws.Range('E10').Value = '  -1.54%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.107'
$ws.Range('E11').Value = '  +2.95%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.156'
$ws.Range('E12').Value = '  -1.35%  '

# Row 13
Set-TextValue $ws.Range('D13') '0.370'
$ws.Range('E13').Value = '  +7.17%  '

# Row 14
Set-TextValue $ws.Range('D14') '3.077.31'
$ws.Range('E14').Value = '  +0.93%  '

# Row 15
Set-TextValue $ws.Range('D15') '60.916.69'
$ws.Range('E15').Value = '  +2.60%  '

# Row 16
Set-TextValue $ws.Range('D16') '23.73'
$ws.Range('E16').Value = '  +5.06%  '

# Row 17
Set-TextValue $ws.Range('D17') '0.0000141'
$ws.Range('E17').Value = '  +2.91%  '

# Row 18
Set-TextValue $ws.Range('D18') '2.626.94'
$ws.Range('E18').Value = '  +1.36%  '

# Row 19
Set-TextValue $ws.Range('D19') '4.71'
$ws.Range('E19').Value = '  +3.81%  '

# Row 20
$ws.Range('E20').Value = '  +9.79%  '

# Row 21
Set-TextValue $ws.Range('D21') '349.71'
$ws.Range('E21').Value = '  +4.26%  '

# Row 22
Set-TextValue $ws.Range('D22') '7.08'
$ws.Range('E22').Value = '  +14.13%  '

# Row 23
$ws.Range('E23').Value = '  +0.32%  '

# Row 24
Set-TextValue $ws.Range('D24') '0.517'
$ws.Range('E24').Value = '  +13.00%  '

# Row 25
Set-TextValue $ws.Range('D25') '63.87'
$ws.Range('E25').Value = '  -0.86%  '

# Row 26
$ws.Range('E26').Value = '  +0.48%  '

# Row 27
$ws.Range('E27').Value = '  +0.35%  '

# Row 28
$ws.Range('E28').Value = '  +6.24%  '

# Row 29
Set-TextValue $ws.Range('D29') '0.0₃0801'
$ws.Range('E29').Value = '  +2.15%  '

# Row 30
$ws.Range('E30').Value = '  +11.51%  '

# Row 31
$ws.Range('E31').Value = '  -0.10%  '

# Row 32
$ws.Range('E32').Value = '  +2.98%  '

# Row 33
Set-TextValue $ws.Range('D33') '162.11'
$ws.Range('E33').Value = '  +2.19%  '

# Row 34
$ws.Range('E34').Value = '  +2.50%  '

# Row 35
Set-TextValue $ws.Range('D35') '4.28'
$ws.Range('E35').Value = '  +5.09%  '

# Row 36
Set-TextValue $ws.Range('D36') '0.956'
$ws.Range('E36').Value = '  +8.38%  '

# Row 37
$ws.Range('E37').Value = '  +4.38%  '

# Row 38
$ws.Range('E38').Value = '  +6.91%  '

# Row 39
Set-TextValue $ws.Range('D39') '37.69'
$ws.Range('E39').Value = '  +1.42%  '

# Row 40
$ws.Range('E40').Value = '  -2.18%  '

# Row 41
$ws.Range('E41').Value = '  +3.48%  '

# Row 42
Set-TextValue $ws.Range('D42') '297.84'
$ws.Range('E42').Value = '  +0.96%  '

# Row 43
Set-TextValue $ws.Range('D43') '139.23'
$ws.Range('E43').Value = '  +10.83%  '

# Row 44
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D44') '0.0988'
$ws.Range('E44').Value = '  +0.86%  '

# Row 45
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D45') '0.995'
$ws.Range('E45').Value = '  -0.33%  '

# Row 46
Set-TextValue $ws.Range('D46') '0.607'
$ws.Range('E46').Value = '  +1.72%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.0552'
$ws.Range('E47').Value = '  +2.23%  '

# Row 48
$ws.Range('E48').Value = '  +3.74%  '

# Row 49
Set-TextValue $ws.Range('D49') '10.71'
$ws.Range('E49').Value = '  +0.66%  '

# Row 50
Set-TextValue $ws.Range('D50') '19.71'
$ws.Range('E50').Value = '  +5.87%  '

# Row 51
Set-TextValue $ws.Range('D51') '2.043.04'
$ws.Range('E51').Value = '  +4.50%  '
